# Excel COM-interop script: add "config_actors" sheet + wire it into
# config_netConnections via a new "owner_actor" column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "config_actors" worksheet right before
#    "config_energyAssets" (so final order is: config, config_netNodes,
#    config_netConnections, config_actors, config_energyAssets).
# ---------------------------------------------------------------------------
$wsEnergyAssets = $wb.Worksheets.Item("config_energyAssets")
$wsActors = $wb.Worksheets.Add($wsEnergyAssets)
$wsActors.Name = "config_actors"

# ---------------------------------------------------------------------------
# 2. Populate config_actors headers + data.
# ---------------------------------------------------------------------------
$wsActors.Cells.Item(1, 1).Value = "index"
$wsActors.Cells.Item(1, 2).Value = "actortype"
$wsActors.Cells.Item(1, 3).Value = "type"
$wsActors.Cells.Item(1, 4).Value = "id"
$wsActors.Cells.Item(1, 5).Value = "parent_actor"

$actorRows = @(
    @(0,  "household",      "HOUSEHOLD",      "hh1",  "sup1"),
    @(1,  "household",      "HOUSEHOLD",      "hh2",  "sup1"),
    @(2,  "household",      "HOUSEHOLD",      "hh3",  "sup1"),
    @(3,  "household",      "HOUSEHOLD",      "hh4",  "hol1"),
    @(4,  "household",      "HOUSEHOLD",      "hh5",  "sup1"),
    @(5,  "household",      "HOUSEHOLD",      "hh6",  "sup1"),
    @(6,  "household",      "HOUSEHOLD",      "hh7",  "sup1"),
    @(7,  "household",      "HOUSEHOLD",      "hh8",  "sup1"),
    @(8,  "household",      "HOUSEHOLD",      "hh9",  "sup1"),
    @(9,  "household",      "HOUSEHOLD",      "hh10", "hol1"),
    @(10, "household",      "HOUSEHOLD",      "hh11", "hol1"),
    @(11, "household",      "HOUSEHOLD",      "hh12", "hol1"),
    @(12, "household",      "HOUSEHOLD",      "hh13", "sup1"),
    @(13, "household",      "HOUSEHOLD",      "hh14", "sup1"),
    @(14, "household",      "HOUSEHOLD",      "hh15", "sup1"),
    @(15, "household",      "HOUSEHOLD",      "hh16", "sup1"),
    @(16, "household",      "HOUSEHOLD",      "hh17", "sup1"),
    @(17, "household",      "HOUSEHOLD",      "hh18", "sup1"),
    @(18, "commercial",     "COMMERCIAL",     "com1", "hol1"),
    @(19, "commercial",     "COMMERCIAL",     "com2", "sup2"),
    @(20, "energysupplier", "ENERGYSUPPLIER", "sup1", "nat"),
    @(21, "energysupplier", "ENERGYSUPPLIER", "sup2", "nat"),
    @(22, "holon",          "HOLON",          "hol1", "nat")
)

for ($i = 0; $i -lt $actorRows.Length; $i++) {
    $r = $i + 2
    $row = $actorRows[$i]
    $wsActors.Cells.Item($r, 1).Value = $row[0]
    $wsActors.Cells.Item($r, 2).Value = $row[1]
    $wsActors.Cells.Item($r, 3).Value = $row[2]
    $wsActors.Cells.Item($r, 4).Value = $row[3]
    $wsActors.Cells.Item($r, 5).Value = $row[4]
}

# Column widths matching the authored layout.
$wsActors.Columns.Item(2).ColumnWidth = 19.33203125
$wsActors.Columns.Item(3).ColumnWidth = 19.33203125
$wsActors.Columns.Item(5).ColumnWidth = 14.1640625

# Selection left on the sheet after data entry.
$wsActors.Range("C25").Select()

# ---------------------------------------------------------------------------
# 3. Add "owner_actor" column (I) to config_netConnections, linking each
#    net connection to the actor that owns it.
# ---------------------------------------------------------------------------
$wsNetConnections = $wb.Worksheets.Item("config_netConnections")
$wsNetConnections.Cells.Item(1, 9).Value = "owner_actor"

$ownerActor = @(
    "hh1", "hh2", "hh3", "hh4", "hh5", "hh6", "hh7", "hh8", "hh9",
    "hh10", "hh11", "hh12", "hh13", "hh14", "hh15", "hh16", "hh17", "hh18",
    "com1", "com1", "com1", "com1", "com1",
    "com2", "com2",
    "sup1", "sup1", "sup1"
)

for ($i = 0; $i -lt $ownerActor.Length; $i++) {
    $r = $i + 2
    $wsNetConnections.Cells.Item($r, 9).Value = $ownerActor[$i]
}

# ---------------------------------------------------------------------------
# 4. Restore config_netConnections as the active sheet/selection.
# ---------------------------------------------------------------------------
$wsNetConnections.Activate()
$wsNetConnections.Range("I30").Select()
